# Rename the worksheet from "Sheet1" to "Query Results" so the exported
# IC reference workbook matches the sheet-naming convention used when the
# file is opened/refreshed via Inquire.
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets("Sheet1")
$ws.Name = "Query Results"
